# Update test e-mail addresses: the timestamp embedded in each address
# changes from 20251109_020650 to 20251109_022039. The same e-mail
# strings also appear on the "LoginData" sheet, so update both sheets
# so no worksheet is left referencing the stale address.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("UsuariosRegistro")
$ws1.Range("C2").Value = "juan.perez+20251109_022039@test.com"
$ws1.Range("C3").Value = "maria.gonzalez+20251109_022039@test.com"
$ws1.Range("C4").Value = "carlos.rodriguez+20251109_022039@test.com"
$ws1.Range("C5").Value = "ana.martinez+20251109_022039@test.com"
$ws1.Range("C6").Value = "luis.garcia+20251109_022039@test.com"

$ws2 = $wb.Worksheets.Item("LoginData")
$ws2.Range("A2").Value = "juan.perez+20251109_022039@test.com"
$ws2.Range("A3").Value = "maria.gonzalez+20251109_022039@test.com"
